$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 289: WB02AA6099 / JAZZ / BODY SHOP / WORK IN PROGRESS
$ws.Range("A289").Value = 44793
$ws.Range("B289").Value = "WB02AA6099"
$ws.Range("C289").Value = "JAZZ"
$ws.Range("D289").Value = "BODY SHOP"
$ws.Range("E289").Value = "WORK IN PROGRESS"

# Row 290: KA53MC6895 / HONDA CITY / PMS / WORK DONE DELIVERED / 6948 / P PAY
$ws.Range("A290").Value = 44793
$ws.Range("B290").Value = "KA53MC6895"
$ws.Range("C290").Value = "HONDA CITY"
$ws.Range("D290").Value = "PMS"
$ws.Range("E290").Value = "WORK DONE DELIVERED"
$ws.Range("F290").Value = 6948
$ws.Range("G290").Value = "P PAY"

# Row 291: KA51MA9141 / FIGO / PMS / WORK IN PROGRESS
$ws.Range("A291").Value = 44793
$ws.Range("B291").Value = "KA51MA9141"
$ws.Range("C291").Value = "FIGO"
$ws.Range("D291").Value = "PMS"
$ws.Range("E291").Value = "WORK IN PROGRESS"

# Row 292: KA03MP6863 / FIGO / PMS / WORK DONE DELIVERED / 8983
$ws.Range("A292").Value = 44793
$ws.Range("B292").Value = "KA03MP6863"
$ws.Range("C292").Value = "FIGO"
$ws.Range("D292").Value = "PMS"
$ws.Range("E292").Value = "WORK DONE DELIVERED"
$ws.Range("F292").Value = 8983

# Row 293: KA03NC9110 / TIAGO / PMS / WORK IN PROGRESS
$ws.Range("A293").Value = 44793
$ws.Range("B293").Value = "KA03NC9110"
$ws.Range("C293").Value = "TIAGO"
$ws.Range("D293").Value = "PMS"
$ws.Range("E293").Value = "WORK IN PROGRESS"

# Row 294: DL8CZ9531 / SCORPIO / SCANNING / WORK IN PROGRESS
$ws.Range("A294").Value = 44793
$ws.Range("B294").Value = "DL8CZ9531"
$ws.Range("C294").Value = "SCORPIO"
$ws.Range("D294").Value = "SCANNING"
$ws.Range("E294").Value = "WORK IN PROGRESS"

$ws.Range("E294").Select()
